$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$abstract = @"
Background
Since December 2019, 2019 novel coronavirus pneumonia emerged in Wuhan city and rapidly spread throughout China and even the world.
 We sought to analyse the clinical characteristics and laboratory findings of some cases with 2019 novel coronavirus pneumonia .
Methods
In this retrospective study, we extracted the data on 95 patients with laboratory-confirmed 2019 novel coronavirus pneumonia in Wuhan Xinzhou District People’s Hospital from January 16th to February 25th, 2020. Cases were confirmed by real-time RT-PCR and abnormal radiologic findings.
 Outcomes were followed up until March 2th, 2020.
Results
Higher temperature, blood leukocyte count, neutrophil count, neutrophil percentage, C-reactive protein level, D-dimer level, alanine aminotransferase activity, aspartate aminotransferase activity, α - hydroxybutyrate dehydrogenase activity, lactate dehydrogenase activity and creatine kinase activity were related to severe 2019 novel coronavirus pneumonia and composite endpoint, and so were lower lymphocyte count, lymphocyte percentage and total protein level.
 Age below 40 or above 60 years old, male, higher Creatinine level, and lower platelet count also seemed related to severe 2019 novel coronavirus pneumonia and composite endpoint, however the P values were greater than 0.05, which mean under the same condition studies of larger samples are needed in the future.
Conclusion
Multiple factors were related to severe 2019 novel coronavirus pneumonia and composite endpoint, and more related studies are needed in the future.

"@

$authors17 = @"
[Gemin%Zhang%NULL%0,        Jie%Zhang%945128911@qq.com%2,        Bowen%Wang%NULL%2,        Xionglin%Zhu%NULL%2,        Qiang%Wang%NULL%2,        Shiming%Qiu%NULL%2]
"@

$authors18 = @"
[Gemin%Zhang%NULL%0,        Jie%Zhang%945128911@qq.com%0,        Bowen%Wang%NULL%0,        Xionglin%Zhu%NULL%0,        Qiang%Wang%NULL%0,        Shiming%Qiu%NULL%0]
"@

$ws.Range("D17").Value2 = $abstract
$ws.Range("E17").Value2 = $authors17
$ws.Range("D18").Value2 = $abstract
$ws.Range("E18").Value2 = $authors18
